# Insert a new weekly price record as row 11 on the active sheet.
# This pushes the existing rows 11-54 down to 12-55 (dimension grows
# from A1:R54 to A1:R55), and fills the freshly inserted row 11 with
# the new "Poroto verde" observation (Sin especificar / Región del
# Maule / $/saco 25 kilos) described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 11..54 down one slot, creating a blank row 11.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the new data record.
$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11, 3).Value = "Ñuble"
$ws.Cells.Item(11, 4).Value = 44560
$ws.Cells.Item(11, 5).Value = 16
$ws.Cells.Item(11, 6).Value = 100112031
$ws.Cells.Item(11, 7).Value = "Poroto verde"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 35000
$ws.Cells.Item(11, 12).Value = 36000
$ws.Cells.Item(11, 13).Value = 35500
$ws.Cells.Item(11, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(11, 15).Value = "Región del Maule"
$ws.Cells.Item(11, 16).Value = 1420
$ws.Cells.Item(11, 17).Value = 25
$ws.Cells.Item(11, 18).Value = "Hortaliza"
